$d = $word.ActiveDocument

$replacements = @(
    ,@("2023-05-30 Tuesday", "2023-05-31 Wednesday")
    ,@("34-20=", "1+32=")
    ,@("46-29=", "98-84=")
    ,@("24+64=", "47+46=")
    ,@("29+27=", "27+35=")
    ,@("24+14=", "74-25=")
    ,@("17+63=", "68-28=")
    ,@("83-29=", "6-5=")
    ,@("74+17=", "72-32=")
    ,@("45+24=", "58+14=")
    ,@("75-22=", "85-65=")
    ,@("30+57=", "55-32=")
    ,@("66-53=", "48+43=")
    ,@("71-41=", "72+20=")
    ,@("28+58=", "25+58=")
    ,@("44+15=", "28+31=")
    ,@("63+7=", "17+45=")
    ,@("57+15=", "78-17=")
    ,@("86+1=", "61-3=")
    ,@("77-17=", "64-41=")
    ,@("41+55=", "82-52=")
    ,@("45-11=", "39+15=")
    ,@("23-20=", "5+27=")
    ,@("47+44=", "21+44=")
    ,@("0+64=", "3+9=")
    ,@("64+2=", "94-49=")
    ,@("57+20=", "81-43=")
    ,@("83-55=", "39+58=")
    ,@("98-57=", "77-66=")
    ,@("40-3=", "48-10=")
    ,@("95-89=", "79-61=")
    ,@("78-8=", "74-0=")
    ,@("16-11=", "25+51=")
    ,@("83-0=", "97-0=")
    ,@("68-56=", "85+2=")
    ,@("29+42=", "66-5=")
    ,@("71+21=", "15+25=")
    ,@("73-3=", "5+13=")
    ,@("39-35=", "43+16=")
    ,@("74+9=", "54+7=")
    ,@("34-0=", "2-0=")
    ,@("2+23=", "68-47=")
    ,@("81-59=", "72-2=")
    ,@("22+65=", "91+4=")
    ,@("21+32=", "83-52=")
    ,@("83-58=", "87-55=")
    ,@("74-31=", "8+26=")
    ,@("56+29=", "34+2=")
    ,@("45-4=", "66+21=")
    ,@("52-51=", "6+75=")
    ,@("0+87=", "40+21=")
    ,@("84-1=", "77-21=")
    ,@("47+34=", "10+79=")
    ,@("58-19=", "57-6=")
    ,@("28-8=", "83-64=")
    ,@("17+46=", "25+38=")
    ,@("97-1=", "71+17=")
    ,@("67-18=", "81-58=")
    ,@("4+78=", "35-1=")
    ,@("63+6=", "10+36=")
    ,@("16+13=", "56-54=")
    ,@("82-33=", "71-10=")
    ,@("67+9=", "52+17=")
    ,@("25+47=", "31+67=")
    ,@("69-34=", "72-17=")
    ,@("49+29=", "46-34=")
    ,@("27+47=", "69-60=")
    ,@("71-25=", "1+51=")
    ,@("69-15=", "22-8=")
    ,@("60-31=", "72-55=")
    ,@("17+77=", "42-18=")
    ,@("17+35=", "70+20=")
    ,@("81-76=", "73-16=")
    ,@("19+6=", "94-56=")
    ,@("97-11=", "1+42=")
    ,@("79-47=", "71+20=")
    ,@("26+28=", "16-6=")
    ,@("95-58=", "32+57=")
    ,@("91-2=", "21+78=")
    ,@("48+3=", "26+69=")
    ,@("61-16=", "39+44=")
    ,@("90+1=", "71-43=")
    ,@("14+43=", "91-78=")
    ,@("29-16=", "32-5=")
    ,@("79-5=", "89-40=")
    ,@("25+70=", "44+42=")
    ,@("60-23=", "84-5=")
    ,@("61+27=", "22-1=")
    ,@("0+41=", "50-25=")
    ,@("90-43=", "85+6=")
    ,@("33-2=", "19+63=")
    ,@("93-59=", "77-2=")
    ,@("35-9=", "56+36=")
    ,@("96-25=", "73-43=")
    ,@("44+3=", "10+32=")
    ,@("70+9=", "30+18=")
    ,@("12+27=", "97-82=")
    ,@("5+85=", "73-71=")
    ,@("23+71=", "54+25=")
    ,@("78-14=", "76-21=")
    ,@("3+36=", "50-0=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Replacements applied: $($replacements.Count)"
